$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Insert 6 new blank rows before row 50 ("start_accuracy"), shifting the old
# rows 50-113 down to 56-119.
$ws.Range("A50:E55").EntireRow.Insert(-4121)

# Copy the formatting of row 49 ("peak_effort" row) onto the 6 new rows so
# they inherit the same cell styles (borders/fills/fonts) as the rest of the
# table, matching what Excel does when a row is inserted via the UI.
$ws.Range("A49:E49").Copy()
$ws.Range("A50:E55").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @("prop_discomfort_downs", "Proportion of discomfort ratings that dropped from the previous discomfort rating"),
    @("prop_discomfort_sames", "Proportion of discomfort ratings that held steady with the previous discomfort rating"),
    @("prop_discomfort_ups", "Proportion of discomfort ratings that increased from the previous discomfort rating"),
    @("prop_effort_downs", "Proportion of effort ratings that dropped from the previous effort rating"),
    @("prop_effort_sames", "Proportion of effort ratings that held steady with the previous effort rating"),
    @("prop_effort_ups", "Proportion of effort ratings that increased from the previous effort rating")
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 50 + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Rows.Item($r).RowHeight = 27
}
